$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest years (2008年, 2009年) which currently occupy rows 2:3.
# This shifts all remaining data rows (2010年..2020年) up by two rows.
$ws.Rows("2:3").Delete()

# After the shift, row 13 is empty (the sheet now spans rows 1-12).
# Copy the formatting of the last existing data row (row 12, 2020年) into the
# new row 13 so the new year label cell keeps the same style (s="1").
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Populate the new row with the 2021年 figures.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 6425
$ws.Range("C13").Value = 249
$ws.Range("D13").Value = 5922
$ws.Range("E13").Value = 1351
$ws.Range("F13").Value = 12740
$ws.Range("G13").Value = 5208
$ws.Range("H13").Value = 13124
$ws.Range("I13").Value = 5013
$ws.Range("J13").Value = 2583
$ws.Range("K13").Value = 10435
$ws.Range("L13").Value = 43651
$ws.Range("M13").Value = 30169
$ws.Range("N13").Value = 6104
$ws.Range("O13").Value = 5596
$ws.Range("P13").Value = 14068
$ws.Range("Q13").Value = 4833
$ws.Range("R13").Value = 2938
$ws.Range("S13").Value = 6359
$ws.Range("T13").Value = 2022
$ws.Range("U13").Value = 10840
$ws.Range("V13").Value = 11263
